$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.496.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.88%  "

$ws.Range("D3").Value = "'1.859.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "'311.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'0.4766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.3791"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.08%  "

$ws.Range("D9").Value = "'0.07323"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").Value = "'0.9287"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("D12").Value = "'0.07772"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "'1.890.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "'5.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").Value = "'6.570"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.53%  "

$ws.Range("D16").Value = "'90.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "'0.000008819"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "

$ws.Range("D19").Value = "'1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "'27.519.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.85%  "

$ws.Range("D21").Value = "'14.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "'5.090"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "'10.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "'156.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.18%  "

$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("D27").Value = "'2.008"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").Value = "'115.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").Value = "'4.945"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("D30").Value = "'0.08862"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "'3.326"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").Value = "'1.200"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("D33").Value = "'0.7535"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.87%  "

$ws.Range("D34").Value = "'4.578"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.56%  "

$ws.Range("D35").Value = "'2.713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("D36").Value = "'0.02040"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "

$ws.Range("D37").Value = "'1.121"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("D38").Value = "'0.5593"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.25%  "

$ws.Range("D39").Value = "'0.05313"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("E40").Value = "  +0.24%  "

$ws.Range("D41").Value = "'7.041"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "'8.491"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("D43").Value = "'0.1522"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("D44").Value = "'10.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").Value = "'0.4870"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.82%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").Value = "'103.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("D48").Value = "'1.663"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.35%  "

$ws.Range("D49").Value = "'67.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.40%  "

$ws.Range("D50").Value = "'0.06098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "'0.9101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.99%  "
